$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C7").Value = -12.81329999999999
$ws.Range("A8").Value = -22.35100000000002
$ws.Range("A10").Value = -21.65789999999998
$ws.Range("A12").Value = -21.56840000000001
$ws.Range("C15").Value = -14.41439999999998
$ws.Range("A18").Value = -22.0185
$ws.Range("C18").Value = -11.9176
$ws.Range("C20").Value = -12.0444
$ws.Range("C29").Value = -11.518
$ws.Range("C30").Value = -12.55049999999999
$ws.Range("C31").Value = -12.8469
$ws.Range("A37").Value = -19.8957
$ws.Range("C40").Value = -13.1445
$ws.Range("C50").Value = -13.62829999999999
$ws.Range("A55").Value = -22.3718
$ws.Range("A68").Value = -21.74069999999999
$ws.Range("C68").Value = -11.7207
$ws.Range("C76").Value = -12.1496
$ws.Range("A77").Value = -20.8464
$ws.Range("A78").Value = -20.56129999999998
$ws.Range("A81").Value = -21.7955
$ws.Range("A82").Value = -22.1489
$ws.Range("C87").Value = -13.42629999999999
$ws.Range("C88").Value = -12.80779999999999
$ws.Range("C96").Value = -12.7161
$ws.Range("C98").Value = -12.17679999999999
$ws.Range("C101").Value = -12.76580000000001
$ws.Range("C102").Value = -13.17550000000001
